$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.298.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.803.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5522"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3876"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07583"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.121"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.12"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.208"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.366"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.795.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001068"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06446"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9996"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.986"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.294.43"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.143"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.63"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.400"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.015.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.93"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.121"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1020"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.740"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.667"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2356"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +17.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06297"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.878"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +16.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02321"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.62"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.049"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6398"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").Value = "  +4.45%  "

$ws.Range("E44").Value = "  -2.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5970"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.685"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.85"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.80%  "

$ws.Range("E49").Value = "  +6.68%  "

$ws.Range("E50").Value = "  +4.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06904"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.33%  "
